# Bug fix in landboundaries - correct OneCrossingLandBoundary test
# Updates the "booked hours" formula on the begroting sheet and the
# current cell selection, matching the author's manual edit in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("begroting")

# Update the "hours" formula for the team member row (B22):
# previously 8+34 (=42), now 8+34+34+16 (=92)
$ws.Range("B22").Formula = "=8+34+34+16"

# Recalculate so dependent cells (B23, C23) pick up the new values
$wb.Application.Calculate()

# Move/update the active cell selection to D20, as reflected in the
# saved sheetView state
$ws.Activate()
$ws.Range("D20").Select()
